# Add a new logbook entry (date + comment row, merged comment/initials row)
# mirroring the existing two-row entry pattern used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing two-row log entry (date row + merged
# comment row) down onto the new rows 24:25 so borders/number formats match
# the rest of the log.
$ws.Range("A7:J8").Copy()
$ws.Range("A24:J25").PasteSpecial(-4122) # xlPasteFormats

# Recreate the merge on the new comment/initials row, same as the other
# entries (e.g. B7:J8, B9:J10, ...).
$ws.Range("B24:J25").MergeCells = $true

# New entry: same date as the surrounding entries (2018-03-11 / serial 43170),
# a new comment, and "TS" as the author initials.
$ws.Range("A24").Value2 = 43170
$ws.Range("B24").Value2 = "Moved a bunch of things around and started placing tracks"
$ws.Range("A25").Value2 = "TS"
